# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  -> currently the stock "Office Theme" colour scheme
#                            (used by the notes master)
#   ppt/theme/theme2.xml  -> currently the "Integral" colour scheme
#                            (used by the slide master / main presentation theme)
#
# The authored edit swaps the two themes' contents (Integral <-> Office Theme).
# The only parts of those two theme XML files that actually differ are the
# <a:clrScheme> name attributes and the 12 colour slots (dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink); font/format schemes are identical in both files.
#
# Through the PowerPoint object model the editable surface for this is
# ThemeColorScheme.Colors(i).RGB (PowerPoint resolves the single theme that
# backs the presentation/slide master for this call). We drive every slot to
# the values the target "Office" colour scheme uses, which reproduces the
# effective (visual) swap for the theme actually in force.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

# index -> packed VBA RGB() value (R + G*256 + B*65536) for the "Office" scheme
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
